# Week 16 logging + season sim update (Lions Players Data)
$wb = $excel.ActiveWorkbook
$wsRush = $wb.Worksheets.Item("Rushing")
$wsRecv = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------
# Rushing sheet updates
# ---------------------------------------------------------------

# Row 4 - J.Williams
$wsRush.Range("C4").Value = 65
$wsRush.Range("D4").Value = 33
$wsRush.Range("E4").Value = 15
$wsRush.Range("F4").Value = 14

# Row 8 - C.Reynolds
$wsRush.Range("C8").Value = 17
$wsRush.Range("D8").Value = 12
$wsRush.Range("E8").Value = 8
$wsRush.Range("F8").Value = 2

# Row 9 - A.St. Brown
$wsRush.Range("D9").Value = 2

# Row 10 - K.Raymond
$wsRush.Range("C10").Value = 2

# Row 11 - J.Jefferson
$wsRush.Range("C11").Value = 1
$wsRush.Range("D11").Value = 0
$wsRush.Range("E11").Value = 1

# New row 12 - S.Zylstra (newly logged this week)
$wsRush.Range("A11").Copy()
$wsRush.Range("A12").PasteSpecial(-4122)
$wsRush.Range("A12").Value = 10
$wsRush.Range("B12").Value = "S.Zylstra"
$wsRush.Range("C12").Value = 0
$wsRush.Range("D12").Value = 0
$wsRush.Range("E12").Value = 1
$wsRush.Range("F12").Value = 0

# ---------------------------------------------------------------
# Receiving sheet updates
# ---------------------------------------------------------------

# Row 6 - J.Cabinda
$wsRecv.Range("C6").Value = 5
$wsRecv.Range("D6").Value = 3

# Row 7 - C.Reynolds
$wsRecv.Range("C7").Value = 7
$wsRecv.Range("D7").Value = 6

# Row 9 - A.St. Brown
$wsRecv.Range("C9").Value = 56
$wsRecv.Range("D9").Value = 45
$wsRecv.Range("E9").Value = 9
$wsRecv.Range("F9").Value = 5
$wsRecv.Range("G9").Value = 8
$wsRecv.Range("H9").Value = 5

# Row 10 - K.Raymond
$wsRecv.Range("C10").Value = 48
$wsRecv.Range("D10").Value = 36
$wsRecv.Range("E10").Value = 17
$wsRecv.Range("G10").Value = 6
$wsRecv.Range("H10").Value = 4

# Row 12 - K.Hodge
$wsRecv.Range("E12").Value = 6
$wsRecv.Range("F12").Value = 1

# Row 14 - T.Kennedy
$wsRecv.Range("C14").Value = 4
$wsRecv.Range("D14").Value = 4
$wsRecv.Range("G14").Value = 1
$wsRecv.Range("H14").Value = 1

# Row 16 - B.Wright
$wsRecv.Range("C16").Value = 14
$wsRecv.Range("D16").Value = 9

# Row 17 - J.Reynolds
$wsRecv.Range("C17").Value = 20
$wsRecv.Range("D17").Value = 13
$wsRecv.Range("E17").Value = 11
$wsRecv.Range("F17").Value = 6
$wsRecv.Range("G17").Value = 3

Write-Host "All updates applied"
